$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ AA = 2939; AB = 2326;  AF = 62523 }
    3  = @{ AA = 2855; AB = 2775;  AF = 67021 }
    4  = @{ AA = 2862; AB = 3377;  AF = 69836 }
    5  = @{ AA = 2725; AB = 3744;  AF = 69901 }
    6  = @{ AA = 2807; AB = 4387;  AF = 73500 }
    7  = @{ AA = 2913; AB = 4944;  AF = 75962 }
    8  = @{ AA = 2891; AB = 5542;  AF = 78362 }
    9  = @{ AA = 2831; AB = 6245;  AF = 81446 }
    10 = @{ AA = 3201; AB = 7267;  AF = 87021 }
    11 = @{ AA = 3561; AB = 8266;  AF = 91643 }
    12 = @{ AA = 3815; AB = 9706;  AF = 97269 }
    13 = @{ AA = 4475; AB = 10793; AF = 101673 }
    14 = @{ AA = 4999; AB = 11398; AF = 105043 }
    15 = @{ AA = 5041; AB = 12055; AF = 103618 }
    16 = @{ AA = 5376; AB = 12534; AF = 109015 }
    17 = @{ AA = 5917; AB = 13627; AF = 115388 }
    18 = @{ AA = 6456; AB = 14423; AF = 121319 }
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("AA$row").Value = $values.AA
    $ws.Range("AB$row").Value = $values.AB
    $ws.Range("AF$row").Value = $values.AF
}
